# TestChiKwadrat.xlsx update
#  - New observed counts (column C) on the "OBB" and "Metoda dokładna" sheets
#    (dependent formulas in D:J recompute automatically).
#  - View/selection bookkeeping: which sheet/cell was active when the file
#    was last saved, and the zoom level on "AABB".

$wb = $excel.ActiveWorkbook

$wsAABB = $wb.Worksheets.Item("AABB")
$wsOBB  = $wb.Worksheets.Item("OBB")
$wsMD   = $wb.Worksheets.Item("Metoda dokładna")

# ---- OBB: updated observed frequencies ------------------------------------
$wsOBB.Range("C4").Value  = 7
$wsOBB.Range("C5").Value  = 44
$wsOBB.Range("C7").Value  = 110
$wsOBB.Range("C8").Value  = 270
$wsOBB.Range("C9").Value  = 69
$wsOBB.Range("C10").Value = 7
$wsOBB.Range("C11").Value = 1
$wsOBB.Range("C12").Value = 1

# ---- Metoda dokładna: updated observed frequencies ------------------------
$wsMD.Range("C3").Value  = 0
$wsMD.Range("C4").Value  = 5
$wsMD.Range("C5").Value  = 60
$wsMD.Range("C6").Value  = 262
$wsMD.Range("C7").Value  = 436
$wsMD.Range("C8").Value  = 208
$wsMD.Range("C9").Value  = 33
$wsMD.Range("C10").Value = 3

# ---- View state ------------------------------------------------------------
# AABB: zoom 160% -> 100%, selection C1:C1048576 -> F45
$wsAABB.Activate()
$wsAABB.Range("F45").Select()
$excel.ActiveWindow.Zoom = 100

# Metoda dokładna: selection G40 -> I19
$wsMD.Activate()
$wsMD.Range("I19").Select()

# OBB becomes the active/last-selected sheet, selection J27 -> H22
$wsOBB.Activate()
$wsOBB.Range("H22").Select()
